$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 9; this pushes the existing
# rows 9..79 down to 10..80 (and the sheet dimension grows to A1:R80).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new observation.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44537
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112052
$ws.Range("G9").Value = "Albahaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("N9").Value = "`$/docena de matas"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 1333
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = "Hortaliza"
